{"js": "// Replace the date line and the 25 three-digit x one-digit multiplication\n// equations in the practice worksheet table with the new day/values.\nconst replacements = [\n  [\"2024-12-21 Saturday\", \"2024-12-22 Sunday\"],\n  [\"422\u00d78=3376\", \"497\u00d74=1988\"],\n  [\"848\u00d76=5088\", \"678\u00d74=2712\"],\n  [\"536\u00d78=4288\", \"741\u00d76=4446\"],\n  [\"777\u00d73=2331\", \"289\u00d77=2023\"],\n  [\"989\u00d78=7912\", \"729\u00d76=4374\"],\n  [\"978\u00d76=5868\", \"366\u00d77=2562\"],\n  [\"161\u00d78=1288\", \"535\u00d75=2675\"],\n  [\"639\u00d74=2556\", \"929\u00d73=2787\"],\n  [\"105\u00d78=840\", \"734\u00d73=2202\"],\n  [\"280\u00d78=2240\", \"716\u00d74=2864\"],\n  [\"953\u00d73=2859\", \"294\u00d72=588\"],\n  [\"699\u00d78=5592\", \"960\u00d79=8640\"],\n  [\"523\u00d79=4707\", \"169\u00d79=1521\"],\n  [\"878\u00d79=7902\", \"503\u00d79=4527\"],\n  [\"200\u00d79=1800\", \"673\u00d72=1346\"],\n  [\"506\u00d75=2530\", \"968\u00d78=7744\"],\n  [\"966\u00d74=3864\", \"365\u00d76=2190\"],\n  [\"584\u00d75=2920\", \"264\u00d78=2112\"],\n  [\"753\u00d79=6777\", \"617\u00d75=3085\"],\n  [\"993\u00d74=3972\", \"824\u00d75=4120\"],\n  [\"964\u00d73=2892\", \"659\u00d76=3954\"],\n  [\"534\u00d75=2670\", \"154\u00d72=308\"],\n  [\"413\u00d74=1652\", \"903\u00d78=7224\"],\n  [\"396\u00d74=1584\", \"397\u00d74=1588\"],\n  [\"837\u00d79=7533\", \"809\u00d75=4045\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 three-digit x one-digit multiplication\n# equations in the practice worksheet table to the new values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-12-21 Saturday\", \"2024-12-22 Sunday\"),\n  @(\"422\u00d78=3376\", \"497\u00d74=1988\"),\n  @(\"848\u00d76=5088\", \"678\u00d74=2712\"),\n  @(\"536\u00d78=4288\", \"741\u00d76=4446\"),\n  @(\"777\u00d73=2331\", \"289\u00d77=2023\"),\n  @(\"989\u00d78=7912\", \"729\u00d76=4374\"),\n  @(\"978\u00d76=5868\", \"366\u00d77=2562\"),\n  @(\"161\u00d78=1288\", \"535\u00d75=2675\"),\n  @(\"639\u00d74=2556\", \"929\u00d73=2787\"),\n  @(\"105\u00d78=840\", \"734\u00d73=2202\"),\n  @(\"280\u00d78=2240\", \"716\u00d74=2864\"),\n  @(\"953\u00d73=2859\", \"294\u00d72=588\"),\n  @(\"699\u00d78=5592\", \"960\u00d79=8640\"),\n  @(\"523\u00d79=4707\", \"169\u00d79=1521\"),\n  @(\"878\u00d79=7902\", \"503\u00d79=4527\"),\n  @(\"200\u00d79=1800\", \"673\u00d72=1346\"),\n  @(\"506\u00d75=2530\", \"968\u00d78=7744\"),\n  @(\"966\u00d74=3864\", \"365\u00d76=2190\"),\n  @(\"584\u00d75=2920\", \"264\u00d78=2112\"),\n  @(\"753\u00d79=6777\", \"617\u00d75=3085\"),\n  @(\"993\u00d74=3972\", \"824\u00d75=4120\"),\n  @(\"964\u00d73=2892\", \"659\u00d76=3954\"),\n  @(\"534\u00d75=2670\", \"154\u00d72=308\"),\n  @(\"413\u00d74=1652\", \"903\u00d78=7224\"),\n  @(\"396\u00d74=1584\", \"397\u00d74=1588\"),\n  @(\"837\u00d79=7533\", \"809\u00d75=4045\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n"}
